$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62: "Results for this location…" -> "Results for this location" (drop the ellipsis)
$ws.Range("A62").Value = "Results for this location"

# Row 65 (previously a blank styled row) now holds a new English/Somali pair
$ws.Range("A65").Value = "What to Expect at This Location"
$ws.Range("B65").Value = "Waxa Laga filanayo Goobtaan"

# Row 66 is brand new - another English/Somali pair
$ws.Range("A66").Value = "Getting results for your location…"
$ws.Range("B66").Value = "Helida natiijooyinka goobtaada…"

# Match the formatting used by the rest of the table (no wrapped text)
$ws.Range("A66:B66").WrapText = $false

# Leave the selection on the newly extended block, like the authored edit did
$ws.Range("A61:B66").Select() | Out-Null
